$d = $word.ActiveDocument

# --- 1. Merge the "Napisac program ... Zp ... Zp ... pierwszych)" sentence,
#        which was split across several runs with proofErr (spell-check)
#        wrappers around "Zp", into a single plain run.
$old1 = "Napisać program, który wyszukuje wszystkie generatory ciała Zp oraz program, który wylicza odwrotność multiplikatywną w Zp oraz Zn (pamiętaj o obsłudze przypadków, gdy odwrotność nie istnieje; n jest iloczynem dwóch liczb pierwszych)"
$d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $old1, 2) | Out-Null

# --- 2. Remove the lone trailing-space run after the first "Rozwiazanie:"
#        heading (the paragraph now ends right after the colon).
$rng = $d.Content
$rng.Find.Execute("Rozwiązanie:") | Out-Null
$afterColon = $d.Range($rng.End, $rng.End + 1)
if ($afterColon.Text -eq " ") {
    $afterColon.Delete()
}

# --- 3. Merge the "Stworzylem klase ExtendedModular ... PolynomialModular ..."
#        sentence (proofErr wrapped class names) into a single run.
$old3 = "Stworzyłem klasę ExtendedModular, która reprezentuje element ciała rozszerzonego. Po tej klasie dziedziczy PolynomialModular, która jest reprezentacją wielomianową elementu ciała rozszerzonego."
$d.Content.Find.Execute($old3, $false, $false, $false, $false, $false, $true, 1, $false, $old3, 2) | Out-Null

# --- 4. Merge "oraz metody wygenerowania tych reprezentacji..." (proofErr
#        wrapped "ch") into a single run.
$old4 = "oraz metody wygenerowania tych reprezentacji. Ile jest reprezentacji, czym się różnią a w czym są podobne? Wygeneruj przynajmniej dwie reprezentacje elementów ciała"
$d.Content.Find.Execute($old4, $false, $false, $false, $false, $false, $true, 1, $false, $old4, 2) | Out-Null

# --- 5. Merge "z logarytmu Zecha:" (proofErr wrapped "Zecha") into one run.
$old5 = "z logarytmu Zecha:"
$d.Content.Find.Execute($old5, $false, $false, $false, $false, $false, $true, 1, $false, $old5, 2) | Out-Null

# --- 6. Merge "Do powyzszej zaleznosci nie da sie zastosowac logarytmu Zecha,
#        nalezy ja wpierw odpowiednio przeksztalcic:" (proofErr wrapped
#        "Zecha") into one run.
$old6 = "Do powyższej zależności nie da się zastosować logarytmu Zecha, należy ją wpierw odpowiednio przekształcić:"
$d.Content.Find.Execute($old6, $false, $false, $false, $false, $false, $true, 1, $false, $old6, 2) | Out-Null
